$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new "Skip" worksheet right after Sheet1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Skip"

# Header row
$ws2.Range("A1").Value = "ID"
$ws2.Range("B1").Value = "definition"
$ws2.Range("C1").Value = "word"
$ws2.Range("D1").Value = "audio_desc"
$ws2.Range("E1").Value = "audio_fn"
$ws2.Range("F1").Value = "theme"
$ws2.Range("G1").Value = "secondary_theme"

# Data rows (same records as Sheet1, offset down by one row for the header)
$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = "tree"
$ws2.Range("C2").Value = "træ"
$ws2.Range("D2").Value = "Aidan Pine"
$ws2.Range("E2").Value = "tree.mp3"
$ws2.Range("F2").Value = "plants"
$ws2.Range("G2").Value = "noun"

$ws2.Range("A3").Value = 2
$ws2.Range("B3").Value = "word"
$ws2.Range("C3").Value = "ord"
$ws2.Range("D3").Value = "Aidan Pine"
$ws2.Range("E3").Value = "ord.mp3"
$ws2.Range("F3").Value = "abstract"
$ws2.Range("G3").Value = "noun"

$ws2.Range("A4").Value = 3
$ws2.Range("B4").Value = "hello"
$ws2.Range("C4").Value = "hej"
$ws2.Range("D4").Value = "Aidan Pine"
$ws2.Range("E4").Value = "hej.mp3"
$ws2.Range("F4").Value = "greetings"
$ws2.Range("G4").Value = "interjection"

$ws2.Range("A5").Value = 4
$ws2.Range("B5").Value = "goodbye"
$ws2.Range("C5").Value = "farvel"
$ws2.Range("F5").Value = "greetings"
$ws2.Range("G5").Value = "interjection"

# Defined name scoped to the new "Skip" sheet
$ws2.Names.Add("data_1", "=Skip!`$A`$2:`$G`$5")

# Selection state: Skip!B8 selected, then Sheet1!D16 re-selected so Sheet1 stays the active tab
[void]$ws2.Range("B8").Select()
[void]$ws1.Range("D16").Select()
